# Added pdf of schedule
# Remove the 5 festival rows that now have a schedule PDF tracked elsewhere,
# then correct a handful of stats for the remaining rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-to-top so earlier row numbers stay valid while deleting.
$ws.Rows(23).Delete()   # IndyFringe Theatre
$ws.Rows(20).Delete()   # Indianapolis Motor Speedway
$ws.Rows(14).Delete()   # Hogan Farms Pumpkin Patch & Corn Maze
$ws.Rows(7).Delete()    # Conner Prairie
$ws.Rows(3).Delete()    # Arts for Lawrence

# After the deletions, apply the remaining value corrections.
$ws.Range("A3").Value = 11        # Avon Community Heritage Festival
$ws.Range("E17").Value = 14982    # Indianapolis Zoo
$ws.Range("A22").Value = 18       # MasterWorks Festival
$ws.Range("A27").Value = 10       # Spirit & Place Festival
$ws.Range("A33").Value = 17       # Waterman's Family Farm
$ws.Range("D33").Value = 4.4      # Waterman's Family Farm
$ws.Range("E33").Value = 669      # Waterman's Family Farm
